$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.056.83'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.78%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.463.77'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +1.35%  '

$ws.Range('E4').Value = '  +0.45%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '485.88'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.21'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +9.76%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.57%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.508'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.68%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.493.69'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.83%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.79'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.61%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0967'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.55%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.331'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.21%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.123'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.35%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.902.59'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.80%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '56.246.16'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.21%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.05'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +6.19%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000135'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.52%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.487.00'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.99%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.51'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +6.84%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.07'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.39%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '317.83'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.996'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.77'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +6.55%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.75'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.68%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.412'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.12%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.23%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.162'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.57%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.595.56'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.78%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.70'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.24%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0787'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +9.00%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +0.30%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.57'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.24%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.19'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.50'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.25%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.19'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.96%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.37%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.73'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.97%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.863'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +6.75%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.99'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.78%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.51'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.33%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.75%  '

$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0555'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.18%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.607'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.53%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.32'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.15%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.76'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +12.09%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '259.18'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +9.89%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0920'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.90%  '

$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.20'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.13%  '

$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0226'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.50'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.04%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.878.81'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.17%  '
